$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Card Code" column (I) and "LLINE" column (J) were in the wrong
# order, which was causing the "Possibly not at 100%" QC reports to be
# misread. Move column I (Card Code) so it lands after column J (LLINE),
# i.e. swap the two columns - this carries the header text along with
# the column's own width/formatting, same as dragging the column in the
# Excel UI.
$ws.Columns.Item(9).Cut() | Out-Null
$ws.Columns.Item(11).Insert() | Out-Null
